$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("NOVEMBER 21")
$ws2 = $wb.Worksheets.Item("DECEMBER 21")

# ---------------------------------------------------------------------------
# NOVEMBER 21 sheet
# ---------------------------------------------------------------------------

# Tenant name correction (row 20): NAMELESS -> NICHOLAS
$ws1.Range("A20").Value = "NICHOLAS"

# Payment entered for row 20
$ws1.Range("G20").Value = 5000

# Extra arrears / deposit entries
$ws1.Range("I36").Value = 5500
$ws1.Range("I47").Value = 3000
$ws1.Range("I48").ClearContents()

# New calculation cell K46
$ws1.Range("K46").Formula = "=G47-E47"

# View state for NOVEMBER 21
$ws1.Activate()
$win1 = $excel.ActiveWindow
$win1.ScrollRow = 46
$win1.ScrollColumn = 1
$ws1.Range("H77").Select()

# ---------------------------------------------------------------------------
# DECEMBER 21 sheet
# ---------------------------------------------------------------------------

$ws2.Range("G9").Value = 4000
$ws2.Range("G12").Value = 4500
$ws2.Range("G13").Value = 2400
$ws2.Range("G14").Value = 4500
$ws2.Range("G19").Value = 5500
$ws2.Range("G21").Value = 5500
$ws2.Range("G23").Value = 5000
$ws2.Range("G24").Value = 5500
$ws2.Range("G25").Value = 4000
$ws2.Range("G26").Value = 4500
$ws2.Range("G27").Value = 4500

$ws2.Rows.Item(28).RowHeight = 15
$ws2.Range("G28").Value = 5000

$ws2.Range("G32").Value = 4500
$ws2.Range("G33").Value = 4500

# Row 35: tenant changed from AUGUSTINE WACHIRA to JOHN NGURE, with new payments
$ws2.Range("A35").Value = "JOHN NGURE"
$ws2.Range("C35").Value = 5000
$ws2.Range("G35").Value = 15000

$ws2.Range("G40").Value = 4000
$ws2.Range("G42").Value = 4500
$ws2.Range("G43").Value = 4500
$ws2.Range("G44").Value = 4500

# Row 45: payment updated
$ws2.Range("G45").Value = 6000

$ws2.Range("G47").Value = 4500
$ws2.Range("G48").Value = 4500

# Row 59: commission formula replaced with a literal number
$ws2.Range("B59").Formula = "=5000"

# Row 63: totals line filled in, copying the formatting from the neighbouring
# cells (B63 -> F63, C63 -> G63) before writing the new values
$ws2.Range("B63").Copy()
$ws2.Range("F63").PasteSpecial(-4122)
$ws2.Range("C63").Copy()
$ws2.Range("G63").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A63").Value = "PAID ON 10/12"
$ws2.Range("C63").Value = 219695
$ws2.Range("E63").Value = "PAID ON 10/12"
$ws2.Range("G63").Value = 219695

# Row 49: a stray space was typed into G49, turning H49 into a #VALUE! error
$ws2.Range("G49").Value = "                                                                                        "

# View state for DECEMBER 21 (active / selected sheet)
$ws2.Activate()
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 35
$win2.ScrollColumn = 1
$ws2.Range("G45").Select()

# ---------------------------------------------------------------------------
# Workbook-level window state
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.WindowState = -4140
